$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D (Price) and E (Volume) columns so numeric-looking
# strings (e.g. "545.74", "6.61") are preserved as text instead of being
# auto-converted to floating point numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.370.33"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.112.06"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "545.74"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.86"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +5.57%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.105.52"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.28%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.61"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.42%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.462"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000228"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +7.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.17"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.620.99"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.400.88"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.71%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.113.54"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.71"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "486.74"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.50"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.707"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.18"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.75"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.46"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.92%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.25"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.01%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.54"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.92%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.41"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "57.85"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "502.24"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.43"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +7.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.08"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.289.54"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +8.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0407"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0806"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.97%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.76"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +8.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.16"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.258"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.42%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "124.32"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.82%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.30"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.20%  "
$ws.Range("B49").Value = "PEPE"
$ws.Range("C49").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₃0540"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +10.67%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.44"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.39%  "
